# Auto-generated: update cached market-price columns (H-N) on leve profit sheets
# Mirrors a scheduled market-data refresh (Universalis-style snapshot) across the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 9946.388999999999
$ws.Range("I98").Value = 11810.267
$ws.Range("J98").Value = 627
$ws.Range("K98").Value = 11810.267
$ws.Range("L98").Value = 627
$ws.Range("M98").Value = -10312.267
$ws.Range("N98").Value = -3623
$ws.Range("H106").Value = 11441.538
$ws.Range("I106").Value = 12748.546
$ws.Range("J106").Value = 4253
$ws.Range("K106").Value = 12748.546
$ws.Range("L106").Value = 4253
$ws.Range("M106").Value = -12117.546
$ws.Range("N106").Value = -5515
$ws.Range("H122").Value = 9946.388999999999
$ws.Range("I122").Value = 11810.267
$ws.Range("J122").Value = 627
$ws.Range("K122").Value = 35430.801
$ws.Range("L122").Value = 1881
$ws.Range("M122").Value = -32980.801
$ws.Range("N122").Value = -6781
$ws.Range("H137").Value = 1309.4546
$ws.Range("I137").Value = 916.2941
$ws.Range("K137").Value = 2748.8823
$ws.Range("M137").Value = -198.8822999999998
$ws.Range("H138").Value = 1454.6465
$ws.Range("I138").Value = 707.56525
$ws.Range("J138").Value = 1680.7368
$ws.Range("K138").Value = 2122.69575
$ws.Range("L138").Value = 5042.2104
$ws.Range("M138").Value = 3017.30425
$ws.Range("N138").Value = -15322.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3898.029
$ws.Range("I32").Value = 3403.9841
$ws.Range("K32").Value = 3403.9841
$ws.Range("M32").Value = -3116.9841
$ws.Range("H74").Value = 1480.8462
$ws.Range("I74").Value = 809.4091
$ws.Range("J74").Value = 2349.7646
$ws.Range("K74").Value = 809.4091
$ws.Range("L74").Value = 2349.7646
$ws.Range("M74").Value = 64.59090000000003
$ws.Range("N74").Value = -4097.7646
$ws.Range("H77").Value = 1480.8462
$ws.Range("I77").Value = 809.4091
$ws.Range("J77").Value = 2349.7646
$ws.Range("K77").Value = 4047.0455
$ws.Range("L77").Value = 11748.823
$ws.Range("M77").Value = 320.9545000000003
$ws.Range("N77").Value = -20484.823
$ws.Range("H132").Value = 1461.5103
$ws.Range("I132").Value = 1169.3529
$ws.Range("K132").Value = 3508.0587
$ws.Range("M132").Value = -978.0587000000005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 13889732
$ws.Range("I94").Value = 17857864
$ws.Range("K94").Value = 17857864
$ws.Range("M94").Value = -17857413
$ws.Range("H107").Value = 2067
$ws.Range("I107").Value = 1749.1
$ws.Range("J107").Value = 3656.5
$ws.Range("K107").Value = 1749.1
$ws.Range("L107").Value = 3656.5
$ws.Range("M107").Value = 170.9000000000001
$ws.Range("N107").Value = -7496.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1731.1818
$ws.Range("I31").Value = 1550.5
$ws.Range("J31").Value = 1948
$ws.Range("K31").Value = 1550.5
$ws.Range("L31").Value = 1948
$ws.Range("M31").Value = -1255.5
$ws.Range("N31").Value = -2538
$ws.Range("H34").Value = 1731.1818
$ws.Range("I34").Value = 1550.5
$ws.Range("J34").Value = 1948
$ws.Range("K34").Value = 1550.5
$ws.Range("L34").Value = 1948
$ws.Range("M34").Value = -1348.5
$ws.Range("N34").Value = -2352
$ws.Range("H58").Value = 923.1667
$ws.Range("I58").Value = 832.7222
$ws.Range("J58").Value = 1194.5
$ws.Range("K58").Value = 832.7222
$ws.Range("L58").Value = 1194.5
$ws.Range("M58").Value = -629.7222
$ws.Range("N58").Value = -1600.5
$ws.Range("H132").Value = 1225.983
$ws.Range("I132").Value = 904.5192
$ws.Range("J132").Value = 3614
$ws.Range("K132").Value = 2713.5576
$ws.Range("L132").Value = 10842
$ws.Range("M132").Value = -183.5576000000001
$ws.Range("N132").Value = -15902
$ws.Range("H136").Value = 923.1667
$ws.Range("I136").Value = 832.7222
$ws.Range("J136").Value = 1194.5
$ws.Range("K136").Value = 2498.1666
$ws.Range("L136").Value = 3583.5
$ws.Range("M136").Value = 51.83339999999998
$ws.Range("N136").Value = -8683.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1423
$ws.Range("I5").Value = 1624.238
$ws.Range("J5").Value = 819.2857
$ws.Range("K5").Value = 4872.714
$ws.Range("L5").Value = 2457.8571
$ws.Range("M5").Value = -4760.714
$ws.Range("N5").Value = -2681.8571
$ws.Range("H17").Value = 1159.091
$ws.Range("J17").Value = 687.5
$ws.Range("L17").Value = 2062.5
$ws.Range("N17").Value = -2400.5
$ws.Range("H34").Value = 2285.7144
$ws.Range("I34").Value = 1750
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 5250
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -5166
$ws.Range("N34").Value = -9168
$ws.Range("H39").Value = 2574.625
$ws.Range("J39").Value = 2291.8462
$ws.Range("L39").Value = 6875.5386
$ws.Range("N39").Value = -7463.5386
$ws.Range("H55").Value = 3498.3333
$ws.Range("J55").Value = 3498.3333
$ws.Range("L55").Value = 10494.9999
$ws.Range("N55").Value = -10848.9999
$ws.Range("H122").Value = 769.94446
$ws.Range("I122").Value = 480.36365
$ws.Range("J122").Value = 1225
$ws.Range("K122").Value = 4323.27285
$ws.Range("L122").Value = 11025
$ws.Range("M122").Value = -1873.27285
$ws.Range("N122").Value = -15925
$ws.Range("H126").Value = 5207.409
$ws.Range("I126").Value = 2507.5
$ws.Range("J126").Value = 5807.3887
$ws.Range("K126").Value = 7522.5
$ws.Range("L126").Value = 17422.1661
$ws.Range("M126").Value = -2582.5
$ws.Range("N126").Value = -27302.1661
$ws.Range("H131").Value = 29415606
$ws.Range("I131").Value = 83333560
$ws.Range("J131").Value = 5813.364
$ws.Range("K131").Value = 250000680
$ws.Range("L131").Value = 17440.092
$ws.Range("M131").Value = -249995640
$ws.Range("N131").Value = -27520.092
$ws.Range("H133").Value = 2474.3333
$ws.Range("I133").Value = 1754
$ws.Range("J133").Value = 3374.75
$ws.Range("K133").Value = 5262
$ws.Range("L133").Value = 10124.25
$ws.Range("M133").Value = -202
$ws.Range("N133").Value = -20244.25
$ws.Range("H134").Value = 4138.154
$ws.Range("I134").Value = 2182.5
$ws.Range("J134").Value = 5007.3335
$ws.Range("K134").Value = 6547.5
$ws.Range("L134").Value = 15022.0005
$ws.Range("M134").Value = -1477.5
$ws.Range("N134").Value = -25162.0005
$ws.Range("H135").Value = 1423
$ws.Range("I135").Value = 1624.238
$ws.Range("J135").Value = 819.2857
$ws.Range("K135").Value = 14618.142
$ws.Range("L135").Value = 7373.571300000001
$ws.Range("M135").Value = -12083.142
$ws.Range("N135").Value = -12443.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2192.3572
$ws.Range("I122").Value = 1874.25
$ws.Range("K122").Value = 5622.75
$ws.Range("M122").Value = -3172.75
$ws.Range("H132").Value = 1767.8049
$ws.Range("I132").Value = 1166.6072
$ws.Range("J132").Value = 3062.6924
$ws.Range("K132").Value = 3499.8216
$ws.Range("L132").Value = 9188.0772
$ws.Range("M132").Value = -969.8215999999998
$ws.Range("N132").Value = -14248.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 880.8333
$ws.Range("I22").Value = 774
$ws.Range("K22").Value = 774
$ws.Range("M22").Value = -479
$ws.Range("H27").Value = 880.8333
$ws.Range("I27").Value = 774
$ws.Range("K27").Value = 774
$ws.Range("M27").Value = -667
$ws.Range("H100").Value = 825.5454999999999
$ws.Range("I100").Value = 557
$ws.Range("K100").Value = 557
$ws.Range("M100").Value = -16
$ws.Range("H136").Value = 1428.1052
$ws.Range("I136").Value = 1324
$ws.Range("J136").Value = 1719.6
$ws.Range("K136").Value = 3972
$ws.Range("L136").Value = 5158.799999999999
$ws.Range("M136").Value = -1422
$ws.Range("N136").Value = -10258.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17858576
$ws.Range("I122").Value = 17858576
$ws.Range("K122").Value = 53575728
$ws.Range("M122").Value = -53573278
$ws.Range("H132").Value = 2621.1667
$ws.Range("I132").Value = 2618.25
$ws.Range("J132").Value = 2644.5
$ws.Range("K132").Value = 7854.75
$ws.Range("L132").Value = 7933.5
$ws.Range("M132").Value = -5324.75
$ws.Range("N132").Value = -12993.5
